$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching original inlineStr cells)
$ws.Range("D2").Value = '67.193.73'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '3.944.73'
$ws.Range("E3").Value = '  +4.04%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '471.30'
$ws.Range("E5").Value = '  +8.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.86'
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.734'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +6.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000338'
$ws.Range("E11").Value = '  +7.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.37'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '4.576.34'
$ws.Range("E13").Value = '  +4.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.36'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.31'
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '3.934.76'
$ws.Range("E16").Value = '  +4.76%  '
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.87'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").Value = '67.486.27'
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.22'
$ws.Range("E21").Value = '  +7.57%  '
$ws.Range("E22").Value = '  +5.10%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.83'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("E25").Value = '  +8.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '39.05'
$ws.Range("E26").Value = '  +5.96%  '
$ws.Range("E27").Value = '  +3.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '722.51'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.56'
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.93'
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.89'
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = '0.0₃0790'
$ws.Range("E37").Value = '  +16.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.35'
$ws.Range("E38").Value = '  -5.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0477'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("E40").Value = '  +4.56%  '
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.52'
$ws.Range("E42").Value = '  +5.21%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.336'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.55'
$ws.Range("E45").Value = '  -7.98%  '
$ws.Range("E46").Value = '  +5.36%  '
$ws.Range("E47").Value = '  +3.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '146.42'
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.12'
$ws.Range("E49").Value = '  -5.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.87'
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.78'
$ws.Range("E51").Value = '  +3.51%  '
